$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.619848132133484
$ws.Range("B1").Value = 4.322647571563721
$ws.Range("C1").Value = 3.527450323104858
$ws.Range("D1").Value = 1.788435697555542
$ws.Range("E1").Value = 1.036090970039368
